# Update "want to go" counts (column F) on the "展览" and "全部类型"
# sheets (and the 3 rows on "本地生活"), mirroring a scheduled refresh
# of the scraped event data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 14671
$ws.Range("F6").Value  = 1417
$ws.Range("F7").Value  = 5973
$ws.Range("F13").Value = 1574
$ws.Range("F15").Value = 2136
$ws.Range("F16").Value = 1247
$ws.Range("F17").Value = 1881
$ws.Range("F20").Value = 2309
$ws.Range("F21").Value = 585
$ws.Range("F22").Value = 842
$ws.Range("F23").Value = 3428
$ws.Range("F26").Value = 2499
$ws.Range("F30").Value = 1847
$ws.Range("F32").Value = 1474
$ws.Range("F35").Value = 5061
$ws.Range("F36").Value = 4983
$ws.Range("F39").Value = 695
$ws.Range("F49").Value = 310

# --- Sheet "本地生活" (local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7759
$ws.Range("F3").Value = 269
$ws.Range("F4").Value = 947

# --- Sheet "全部类型" (all types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 7759
$ws.Range("F4").Value  = 269
$ws.Range("F5").Value  = 947
$ws.Range("F8").Value  = 14671
$ws.Range("F10").Value = 1417
$ws.Range("F11").Value = 5973
$ws.Range("F16").Value = 1574
$ws.Range("F19").Value = 842
$ws.Range("F20").Value = 3428
$ws.Range("F22").Value = 2499
$ws.Range("F25").Value = 1847
$ws.Range("F32").Value = 1474
$ws.Range("F35").Value = 5061
$ws.Range("F36").Value = 4983
$ws.Range("F38").Value = 695
$ws.Range("F46").Value = 310
